# Apply updated cryptocurrency price/volume data to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = "<price text>"; E = "<volume/% text>" }
# A leading apostrophe is used for D values that look like genuine numbers
# ("244.80", "0.9989", ...) so Excel stores them as text, exactly like the
# original inline-string cells (values such as "29.610.32" already contain
# two dots so Excel treats them as text automatically).
$updates = @(
    @{ Row = 2; D = '29.610.32'; E = '  +2.41%  ' }
    @{ Row = 3; D = '1.858.72'; E = '  +1.49%  ' }
    @{ Row = 4; D = '''0.9989'; E = '  -0.08%  ' }
    @{ Row = 5; D = '''244.80'; E = '  +0.21%  ' }
    @{ Row = 6; D = '''0.6942'; E = '  +0.79%  ' }
    @{ Row = 7; D = '''0.9997'; E = '  -0.03%  ' }
    @{ Row = 8; D = '''0.07697'; E = '  +0.49%  ' }
    @{ Row = 9; D = '''0.3059'; E = '  +0.16%  ' }
    @{ Row = 10; D = '''23.72'; E = '  +0.84%  ' }
    @{ Row = 11; D = '''0.07774'; E = '  -0.53%  ' }
    @{ Row = 12; D = '''5.144'; E = '  +1.25%  ' }
    @{ Row = 13; D = '1.850.94'; E = '  +1.00%  ' }
    @{ Row = 14; D = '''91.52'; E = '  +1.07%  ' }
    @{ Row = 15; D = '''0.6926'; E = '  +2.21%  ' }
    @{ Row = 16; D = '''6.570'; E = '  +1.95%  ' }
    @{ Row = 17; D = '29.585.04'; E = '  +2.32%  ' }
    @{ Row = 18; D = '''0.000008292'; E = '  +0.14%  ' }
    @{ Row = 19; D = '2.100.96'; E = '  +0.89%  ' }
    @{ Row = 20; D = '''240.09'; E = '  -1.25%  ' }
    @{ Row = 21; D = '''12.78'; E = '  +0.78%  ' }
    @{ Row = 22; D = '''0.9996'; E = '  -0.02%  ' }
    @{ Row = 23; D = '''7.607'; E = '  +2.18%  ' }
    @{ Row = 24; D = '''0.9999'; E = '  +0.00%  ' }
    @{ Row = 25; D = '''0.1501'; E = '  +1.88%  ' }
    @{ Row = 26; D = '''8.929'; E = '  +1.51%  ' }
    @{ Row = 27; D = '''159.76'; E = '  -0.97%  ' }
    @{ Row = 28; E = '  +0.45%  ' }
    @{ Row = 29; E = '  -1.20%  ' }
    @{ Row = 30; D = '''4.252'; E = '  +0.92%  ' }
    @{ Row = 31; D = '''4.179'; E = '  +1.23%  ' }
    @{ Row = 32; E = '  +2.38%  ' }
    @{ Row = 33; D = '''0.05099'; E = '  -0.44%  ' }
    @{ Row = 34; D = '''0.7716'; E = '  +1.74%  ' }
    @{ Row = 35; D = '''1.896'; E = '  +3.32%  ' }
    @{ Row = 36; D = '''1.153'; E = '  +0.81%  ' }
    @{ Row = 37; D = '''2.683'; E = '  +0.21%  ' }
    @{ Row = 38; D = '1.335.99'; E = '  +8.10%  ' }
    @{ Row = 39; D = '''0.01871'; E = '  +1.45%  ' }
    @{ Row = 40; D = '''2.725'; E = '  +1.56%  ' }
    @{ Row = 41; D = '''0.9718'; E = '  +4.57%  ' }
    @{ Row = 42; D = '''106.55'; E = '  -1.83%  ' }
    @{ Row = 43; D = '''5.811'; E = '  +2.19%  ' }
    @{ Row = 44; D = '''0.9995'; E = '  -0.02%  ' }
    @{ Row = 45; D = '''9.776'; E = '  +2.58%  ' }
    @{ Row = 46; D = '2.001.59'; E = '  +1.09%  ' }
    @{ Row = 47; E = '  +3.40%  ' }
    @{ Row = 48; E = '  +0.93%  ' }
    @{ Row = 49; D = '''1.779'; E = '  +2.44%  ' }
    @{ Row = 50; D = '''63.66'; E = '  -0.95%  ' }
    @{ Row = 51; D = '''6.959'; E = '  +0.82%  ' }
)

foreach ($update in $updates) {
    if ($update.ContainsKey("D")) {
        $ws.Range("D" + $update.Row).Value = $update.D
    }
    if ($update.ContainsKey("E")) {
        $ws.Range("E" + $update.Row).Value = $update.E
    }
}
